# Apply the row reshuffle described by the diff: rows 2-23 (columns A-F) are
# permuted to a new order. Row 1 (headers) and rows 24-26 (totals) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..23, columns A..F (after the reorder)
$data = @(
    @(701,  3, 90, 45, 97, 15),
    @(801,  3, 67, 65, 52, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(101,  9, 30, 15, 60, 15),
    @(401,  9, 48, 67, 75, 45),
    @(1001, 18, 30, 75, 60, 72),
    @(501,  9, 52, 30, 75, 45),
    @(901,  16, 15, 45, 60, 60),
    @(902,  1, 0, 0, 0, 0),
    @(201,  9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(301,  6, 45, 30, 60, 45),
    @(601,  9, 60, 67, 60, 42),
    @(2,    0, 2, 2, 2, 2),
    @(1,    0, 2, 2, 2, 2),
    @(502,  0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(802,  0, 4, 5, 4, 0),
    @(3,    0, 3, 3, 3, 3),
    @(402,  0, 0, 4, 0, 0),
    @(602,  0, 0, 4, 0, 9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
